$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update utilisation amounts (column F) for rows 2-5 to the new value
$ws.Range("F2").Value = 761579.37
$ws.Range("F3").Value = 761579.37
$ws.Range("F4").Value = 761579.37
$ws.Range("F5").Value = 761579.37

# Update the active selection on the sheet
$ws.Range("F2").Select()
